$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Through 2022-08-19"

# Update header text
$ws.Range("B1").Value = "August 2022 (through August 19)"

# Update the data cells
$ws.Range("B2").Value = 13
$ws.Range("R2").Value = 5
$ws.Range("AP2").Value = 2
$ws.Range("AX2").Value = 2
$ws.Range("R3").Value = 5
$ws.Range("B4").Value = 3
$ws.Range("AP4").Value = 6
$ws.Range("R5").Value = 10
$ws.Range("B7").Value = 7
$ws.Range("AH7").Value = 2
$ws.Range("R9").Value = 5
$ws.Range("AX9").Value = 8
$ws.Range("Z13").Value = 2
$ws.Range("B15").Value = 6
$ws.Range("R15").Value = 3
$ws.Range("Z15").Value = 2
$ws.Range("R16").Value = 3
$ws.Range("AH16").Value = 1
$ws.Range("J17").Value = 2
$ws.Range("J33").Value = 2
$ws.Range("R45").Value = 3
